$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "latitude"
$ws.Range("B3").Value = "longitude"
$ws.Range("B4").Value = "latitude"
$ws.Range("B5").Value = "longitude"

$validation = $ws.Range("B1:B1048576").Validation
$validation.Delete()
$validation.Add(3, 1, 1, '"latitude,longitude,geography,temporal,genomic"')
$validation.ErrorTitle = "Invalid category type"
$validation.ErrorMessage = "The GenEpi DRIVE application currently only supports the following data types: `n`ngeography, temporal, or genomic`n`nplease assign you variable to one of these categories, or, leave it blank"
$validation.ShowInput = $true
$validation.ShowError = $true

$ws.Range("B8").Select() | Out-Null
